# BA.xlsx TanzyWatch update: add Jun_27 and Jun_26 (x2) date columns at the
# front of the rating grid, and add two new analysts (Benchmark, Evercore ISI)
# as new rows at the bottom of the watch list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new data columns right after the analyst-name column (A),
# pushing the existing Jun_17 / Jun_15 / Jun_13 / Jun_10 columns from B:E to E:H.
$ws.Range("B1:D1").EntireColumn.Insert()

# Give the three new columns the same width/format as their neighbours.
$ws.Columns.Item(3).ColumnWidth = 7.14
$ws.Columns.Item(4).ColumnWidth = 7.14
$ws.Columns.Item(5).ColumnWidth = 7.14
$ws.Columns.Item(6).ColumnWidth = 7.14
$ws.Columns.Item(7).ColumnWidth = 7.14
$ws.Columns.Item(8).ColumnWidth = 7.14

# New header row: two new reporting dates (Jun_27, then Jun_26 twice).
$ws.Cells.Item(1, 2).Value = "Jun_27"
$ws.Cells.Item(1, 3).Value = "Jun_26"
$ws.Cells.Item(1, 4).Value = "Jun_26"

# Fill the new columns with "UN" (unchanged) for every existing analyst row,
# matching the placeholder value already used across the rest of the grid.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# Add the two new analyst rows at the bottom of the watch list.
$ws.Cells.Item(28, 1).Value = "Benchmark"
$ws.Cells.Item(28, 2).Value = "UN"
$ws.Cells.Item(28, 3).Value = "UN"
$ws.Cells.Item(28, 4).Value = "UN"

$ws.Cells.Item(29, 1).Value = "Evercore ISI"
$ws.Cells.Item(29, 2).Value = "UN"
$ws.Cells.Item(29, 3).Value = "UN"
$ws.Cells.Item(29, 4).Value = "UN"
